$d = $word.ActiveDocument

# The paragraph ending with "Os dashboards começam a mostrar dados em
# poucos segundos." is the last paragraph in the document.
$p = $d.Paragraphs.Item($d.Paragraphs.Count)

# 1) A lone line break right after the existing sentence.
$r = $p.Range
$r.Collapse(0)
$r.InsertBreak(6)

# 2) A second line break immediately followed by "acesse o link no "
#    (Chr(11) is a manual line break embedded in the same run as the text).
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter([char]11 + "acesse o link no ")

# 3) "git"
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter("git")

# 4) " hub para clonar o projeto:"
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter(" hub para clonar o projeto:")

# 5) A lone line break.
$r = $p.Range
$r.Collapse(0)
$r.InsertBreak(6)

# 6) The repo URL.
$r = $p.Range
$r.Collapse(0)
$r.InsertAfter("https://github.com/abruno36/grpc-observability-stack")

# Underline just the URL we inserted (everything else in the paragraph
# keeps its original formatting).
$fr = $d.Content
$fr.Find.ClearFormatting()
$fr.Find.Execute("https://github.com/abruno36/grpc-observability-stack", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fr.Font.Underline = 1
